$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 trailing rows (28 -> 23 rows); content will be fully rewritten below.
$ws.Range("A24:A28").EntireRow.Delete()

# Rows 1-9 are unchanged by this edit - leave them untouched.

# Clear stray leftover cells from the old layout that have no counterpart
# in the new layout (rows 10-23 are being restructured). Use Clear() (not
# ClearContents()) so the cell element itself is dropped, not left as an
# empty styled stub.
$ws.Cells.Item(13,2).Clear()
$ws.Cells.Item(13,3).Clear()
$ws.Cells.Item(15,2).Clear()
$ws.Cells.Item(15,3).Clear()
$ws.Cells.Item(16,2).Clear()
$ws.Cells.Item(16,3).Clear()
$ws.Cells.Item(22,1).Clear()
$ws.Cells.Item(23,1).Clear()

$ws.Cells.Item(10,1).Value = "Objetivos:"
$ws.Cells.Item(10,2).Value = "5840726 - Cristina Bormio Nunes"
$ws.Cells.Item(10,3).Value = "5840726 - Cristina Bormio Nunes"
$ws.Rows.Item(10).RowHeight = 60
$ws.Cells.Item(11,1).Value = "Objectives:"
$ws.Rows.Item(11).RowHeight = 60
$ws.Cells.Item(12,1).Value = "Programa resumido:"
$ws.Cells.Item(12,2).Value = "6495737 - Durval Rodrigues Junior"
$ws.Cells.Item(12,3).Value = "6495737 - Durval Rodrigues Junior"
# Column B is covered by two overlapping <col> style ranges (style 1 and
# style 2); nudge the newly-created B12 cell onto the "column B" style
# (s=2, matching B10/B14/...) instead of the ambiguous default.
$ws.Cells.Item(12,2).WrapText = $true
$ws.Cells.Item(12,2).VerticalAlignment = -4160
$ws.Cells.Item(12,2).Font.Bold = $false
$ws.Rows.Item(12).RowHeight = 60
$ws.Cells.Item(13,1).Value = "Short syllabus:"
$ws.Rows.Item(13).RowHeight = 60
$ws.Cells.Item(14,1).Value = "Programa:"
$ws.Cells.Item(14,2).Value = "1341653 - Maria José Ramos Sandim"
$ws.Cells.Item(14,3).Value = "1341653 - Maria José Ramos Sandim"
$ws.Rows.Item(14).RowHeight = 120
$ws.Cells.Item(15,1).Value = "Syllabus:"
$ws.Rows.Item(15).RowHeight = 120
$ws.Cells.Item(16,1).Value = "Avaliação:"
$ws.Cells.Item(17,1).Value = "Método:"
$ws.Cells.Item(17,2).Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Cells.Item(17,3).Value = "1643715 - Paulo Atsushi Suzuki"
$ws.Rows.Item(17).RowHeight = 60
$ws.Cells.Item(18,1).Value = "Critério:"
$ws.Cells.Item(18,2).Value = "Aulas expositivas e  exercícios comentados"
$ws.Cells.Item(18,3).Value = "Aulas expositivas e  exercícios comentados"
$ws.Rows.Item(18).RowHeight = 60
$ws.Cells.Item(19,1).Value = "Norma de recuperação:"
$ws.Cells.Item(19,2).Value = "Média final calculada pelas notas de 2 provas (P1 e P2), seguindo os pesos MF=(P1+2*P2)/3, ou seja, peso 1 para a P1 e peso 2 para a P2."
$ws.Cells.Item(19,3).Value = "Média final calculada pelas notas de 2 provas (P1 e P2), seguindo os pesos MF=(P1+2*P2)/3, ou seja, peso 1 para a P1 e peso 2 para a P2."
$ws.Rows.Item(19).RowHeight = 60
$ws.Cells.Item(20,1).Value = "Bibliografia:"
$ws.Cells.Item(20,2).Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Cells.Item(20,3).Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Rows.Item(20).RowHeight = 120
$ws.Cells.Item(21,1).Value = "Requisitos:"
$ws.Cells.Item(22,2).Value = "LOB1052 -  Cálculo III  (Requisito)`n"
$ws.Cells.Item(22,3).Value = "LOB1052 -  Cálculo III  (Requisito)`n"
$ws.Rows.Item(22).RowHeight = 30
$ws.Cells.Item(23,2).Value = "LOB1053 -  Física III  (Requisito)`n"
$ws.Cells.Item(23,3).Value = "LOB1053 -  Física III  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30
